# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C),
# First_Noticeable_Increase_Cumulative_Value (E) and Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

# Values keyed by worksheet name -> row -> column letter -> new value
$updates = @{
    "Step3_DataPts_0.5" = @{
        2 = @{ "C" = 88; "E" = 0.01933178791621492; "G" = 21 }
        3 = @{ "C" = 87; "E" = 0.01571361524422738; "G" = 38 }
        4 = @{ "C" = 87; "E" = 0.01041399391677734; "G" = 23 }
        5 = @{ "C" = 87; "E" = 0.005092500148287457; "G" = 27 }
        6 = @{ "C" = 88; "E" = 0.04062270445707165; "G" = 20 }
    }
    "Step3_DataPts_0.7" = @{
        2 = @{ "C" = 88; "E" = 0.01933178791621492; "G" = 67 }
        3 = @{ "C" = 87; "E" = 0.01571361524422738; "G" = 68 }
        4 = @{ "C" = 87; "E" = 0.01041399391677734; "G" = 71 }
        5 = @{ "C" = 87; "E" = 0.005092500148287457; "G" = 67 }
        6 = @{ "C" = 88; "E" = 0.04062270445707165; "G" = 63 }
    }
    "Step3_DataPts_0.8" = @{
        2 = @{ "C" = 88; "E" = 0.01933178791621492; "G" = 72 }
        3 = @{ "C" = 87; "E" = 0.01571361524422738; "G" = 73 }
        4 = @{ "C" = 87; "E" = 0.01041399391677734; "G" = 73 }
        5 = @{ "C" = 87; "E" = 0.005092500148287457; "G" = 73 }
        6 = @{ "C" = 88; "E" = 0.04062270445707165; "G" = 68 }
    }
    "Step3_DataPts_0.9" = @{
        2 = @{ "C" = 88; "E" = 0.01933178791621492; "G" = 83 }
        3 = @{ "C" = 87; "E" = 0.01571361524422738; "G" = 84 }
        4 = @{ "C" = 87; "E" = 0.01041399391677734; "G" = 84 }
        5 = @{ "C" = 87; "E" = 0.005092500148287457; "G" = 83 }
        6 = @{ "C" = 88; "E" = 0.04062270445707165; "G" = 79 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
        }
    }
}
